$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the style on AB87: it was style 9 (an accidental one-off font tweak)
#    and should match the rest of the "Kод" column (style 6), per
#    "fix new sku type" in the commit message.
# ---------------------------------------------------------------------------
$ws.Range("AB82").Copy()
$ws.Range("AB87").PasteSpecial(-4122, $false, $false, $false)

# ---------------------------------------------------------------------------
# 2. Append the new SKU as row 88 (Сулугуни в рассоле "Вкусвилл").
# ---------------------------------------------------------------------------

# Row-number style (bold Cambria on a thin border) matches the rest of col A.
$ws.Range("A87").Copy()
$ws.Range("A88").PasteSpecial(-4122, $false, $false, $false)
$ws.Range("A88").Value = 86

# SKU name needs the wrapping style used by the other multi-line names.
$ws.Range("B88").Value = 'Сулугуни в рассоле "Вкусвилл", 45%, 0,21/0,35 кг, ф/п'
$ws.Range("B88").WrapText = $true

$ws.Range("C88").Value = 2.7
$ws.Range("D88").Value = "Да"
$ws.Range("E88").Value = "Сулугуни"
$ws.Range("F88").Value = "Соль"
$ws.Range("G88").Value = "Сакко"
$ws.Range("H88").Value = "Foodfest"
$ws.Range("I88").Value = 210
$ws.Range("J88").Value = 6
$ws.Range("K88").Value = 370
$ws.Range("L88").Value = 960
$ws.Range("M88").Value = 50
$ws.Range("N88").Value = "Нет"
$ws.Range("O88").Value = "Ульма"
$ws.Range("Q88").Value = 1300
$ws.Range("R88").Value = 1300
$ws.Range("U88").Value = 90
$ws.Range("V88").Value = 40
$ws.Range("W88").Value = 20
$ws.Range("X88").Value = 20
$ws.Range("Y88").Value = 15
$ws.Range("Z88").Value = 5
$ws.Range("AA88").Value = 5

# "Kод" cell picks up the same (corrected) style as the rest of the column.
$ws.Range("AB82").Copy()
$ws.Range("AB88").PasteSpecial(-4122, $false, $false, $false)
$ws.Range("AB88").Value = "00-00013255"

$ws.Range("AC88").Value = 1300
$ws.Range("AD88").Value = "Нет"

# ---------------------------------------------------------------------------
# 3. Move the cursor/selection the way the author left it (AB89) after
#    entering the new row.
# ---------------------------------------------------------------------------
$ws.Range("AB89").Select() | Out-Null
